# Fix typographic middle-dots ("·", U+00B7) used as decimal points in
# percentage values on the "Table" worksheet, replacing them with normal
# periods (".") so the numbers read correctly, e.g. "161 (36·3%)" -> "161 (36.3%)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

$ws.Range("B2").Value  = "161 (36.3%)"
$ws.Range("B31").Value = "59 (13.3%)"
$ws.Range("B41").Value = "223 (50.3%)"
$ws.Range("B54").Value = "4 (23.5%)"
$ws.Range("B55").Value = "4 (23.5%)"
